# Refresh the Price (D) and Volume(1h) (E) columns of the cryptos list with the
# latest scrape values. Matches the GitHub Actions "Updated cryptos list" job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates: each entry holds the new Price/Volume text for that row.
# A $null column means that column did not change for that coin this run.
$updates = @(
    @{ Row=2; D="66.420.14"; E="  -1.21%  " },
    @{ Row=3; D="3.444.94"; E="  -0.81%  " },
    @{ Row=4; D=$null; E="  +0.08%  " },
    @{ Row=5; D="579.14"; E="  -2.27%  " },
    @{ Row=6; D="175.09"; E="  -1.67%  " },
    @{ Row=7; D=$null; E="  +0.07%  " },
    @{ Row=8; D="0.598"; E="  +1.31%  " },
    @{ Row=9; D="3.443.97"; E="  -0.91%  " },
    @{ Row=10; D="0.133"; E="  -2.87%  " },
    @{ Row=11; D="6.84"; E="  -3.24%  " },
    @{ Row=12; D="0.418"; E="  -3.11%  " },
    @{ Row=13; D="4.044.78"; E="  -0.73%  " },
    @{ Row=14; D="30.73"; E="  -3.62%  " },
    @{ Row=15; D="0.131"; E="  -3.32%  " },
    @{ Row=16; D="66.414.64"; E="  -1.31%  " },
    @{ Row=17; D="0.0000171"; E="  -3.10%  " },
    @{ Row=18; D="3.450.01"; E="  -0.78%  " },
    @{ Row=19; D="5.98"; E="  -4.03%  " },
    @{ Row=20; D="13.79"; E="  -3.14%  " },
    @{ Row=21; D="375.15"; E="  -3.32%  " },
    @{ Row=22; D="7.67"; E="  -2.25%  " },
    @{ Row=23; D="0.999"; E="  +0.09%  " },
    @{ Row=24; D=$null; E="  +0.16%  " },
    @{ Row=25; D="70.69"; E="  -3.74%  " },
    @{ Row=26; D="0.525"; E="  -1.57%  " },
    @{ Row=27; D="0.0000116"; E="  -3.36%  " },
    @{ Row=28; D="9.80"; E="  -5.10%  " },
    @{ Row=29; D="0.172"; E="  -1.80%  " },
    @{ Row=30; D=$null; E="  -0.02%  " },
    @{ Row=31; D="5.82"; E="  -5.21%  " },
    @{ Row=32; D="23.80"; E="  +1.41%  " },
    @{ Row=33; D=$null; E="  -3.65%  " },
    @{ Row=34; D="1.33"; E="  -5.78%  " },
    @{ Row=35; D=$null; E="  -0.08%  " },
    @{ Row=36; D="7.02"; E="  -4.45%  " },
    @{ Row=37; D=$null; E="  -4.91%  " },
    @{ Row=38; D="159.35"; E="  -2.79%  " },
    @{ Row=39; D="0.876"; E="  +0.68%  " },
    @{ Row=40; D="26.99"; E="  +2.59%  " },
    @{ Row=41; D=$null; E="  -5.02%  " },
    @{ Row=42; D="2.61"; E="  -3.52%  " },
    @{ Row=43; D="6.49"; E="  -5.18%  " },
    @{ Row=44; D="4.44"; E="  -3.64%  " },
    @{ Row=45; D="2.692.23"; E="  -4.46%  " },
    @{ Row=46; D="0.0691"; E="  -3.78%  " },
    @{ Row=47; D="25.15"; E="  -5.70%  " },
    @{ Row=48; D="40.18"; E="  -3.15%  " },
    @{ Row=49; D="0.0293"; E="  -1.80%  " },
    @{ Row=50; D="319.92"; E="  -4.74%  " },
    @{ Row=51; D="1.01"; E="  -3.74%  " }
)

# Price values are plain text in this sheet (e.g. "66.420.14", "1.00", "0.598")
# and must stay text - otherwise Excel would coerce numeric-looking strings
# (trailing zeros, thousand-grouped "xx.xxx.xx" prices, etc.) into numbers and
# mangle them. Prefixing with a leading apostrophe forces a text literal for
# anything that parses as a plain number; values that are already unambiguous
# text (the multi-dot big-coin prices) are written as-is.
function Set-TextCell {
    param($Range, [string]$Value)
    if ($Value -match '^[+-]?\d+(\.\d+)?$') {
        $Range.Value = "'" + $Value
    } else {
        $Range.Value = $Value
    }
}

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextCell $ws.Range("D$($u.Row)") $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
